$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("A13").Value = 251984
$ws.Range("C13").Value = 25
$ws.Range("D13").Value = 338.327868852459
$ws.Range("F13").Value = "2025-06-04 07:25:00"
$ws.Range("G13").Value = "2025-06-04 07:25:00"
$ws.Range("H13").Value = "2025-06-04 13:03:19"
$ws.Range("I13").Value = 20638
$ws.Range("L13").Value = 3
$ws.Range("N13").Value = 39874
$ws.Range("P13").Value = 39874
$ws.Range("Q13").Value = "2025-06-10 00:00:00"
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 1

# Row 14
$ws.Range("A14").Value = 251180
$ws.Range("C14").Value = 25
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = "2025-06-04 13:03:19"
$ws.Range("F14").Value = "2025-06-04 13:28:19"
$ws.Range("G14").Value = "2025-06-04 13:28:19"
$ws.Range("H14").Value = "2025-06-04 13:28:19"
$ws.Range("I14").Value = 0
$ws.Range("L14").Value = 4
$ws.Range("N14").Value = "39887 (esterno)"
$ws.Range("P14").Value = 39887
$ws.Range("Q14").Value = "2025-05-20 00:00:00"
$ws.Range("R14").Value = -15.5613387978125
$ws.Range("S14").Value = 7

# Row 15
$ws.Range("A15").Value = 252282
$ws.Range("D15").Value = 44.88524590163934
$ws.Range("E15").Value = "2025-06-04 13:28:19"
$ws.Range("F15").Value = "2025-06-04 13:53:19"
$ws.Range("G15").Value = "2025-06-04 13:53:19"
$ws.Range("H15").Value = "2025-06-04 14:38:12"
$ws.Range("I15").Value = 2738
$ws.Range("L15").Value = 5
$ws.Range("N15").Value = 39885
$ws.Range("P15").Value = 39885
$ws.Range("Q15").Value = "2025-06-09 00:00:00"
$ws.Range("R15").Value = 0
$ws.Range("S15").Value = 1

# Row 16
$ws.Range("A16").Value = 252084
$ws.Range("C16").Value = 35
$ws.Range("D16").Value = 641
$ws.Range("E16").Value = "2025-06-04 14:38:12"
$ws.Range("F16").Value = "2025-06-05 07:13:12"
$ws.Range("G16").Value = "2025-06-05 07:13:12"
$ws.Range("H16").Value = "2025-06-06 09:54:12"
$ws.Range("I16").Value = 39101
$ws.Range("L16").Value = 2
$ws.Range("N16").Value = 39885
$ws.Range("P16").Value = 39885
$ws.Range("Q16").Value = "2025-06-30 00:00:00"
$ws.Range("R16").Value = -1.412647996354166
$ws.Range("S16").Value = 7

# Row 36
$ws.Range("C36").Value = 35
$ws.Range("E36").Value = "2025-06-06 09:54:12"
$ws.Range("F36").Value = "2025-06-06 10:29:12"
$ws.Range("G36").Value = "2025-06-06 10:29:12"
$ws.Range("H36").Value = "2025-06-10 08:31:20"

# Row 57
$ws.Range("E57").Value = "2025-06-10 08:31:20"
$ws.Range("F57").Value = "2025-06-10 08:56:20"
$ws.Range("G57").Value = "2025-06-10 08:56:20"
$ws.Range("H57").Value = "2025-06-11 13:18:00"
$ws.Range("R57").Value = -6.554178051006944

# Row 70
$ws.Range("E70").Value = "2025-06-11 13:18:00"
$ws.Range("F70").Value = "2025-06-11 13:48:00"
$ws.Range("G70").Value = "2025-06-11 13:48:00"
$ws.Range("H70").Value = "2025-06-12 10:08:00"
